$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Highlight the DTR rows with a light-blue fill (#29A3CC).
#    Rows 5-6 and 11-15 (columns A:J) are the two blocks that get colored.
#    Excel's Interior.Color is a BGR long, so 0x29A3CC (RRGGBB) -> 0xCCA329.
# ---------------------------------------------------------------------------
$blue = 13411113   # RGB(0x29,0xA3,0xCC) packed as BGR long

$ws.Range("A5:J6").Interior.Color = $blue
$ws.Range("A11:J15").Interior.Color = $blue

# ---------------------------------------------------------------------------
# 2) B19 should hold the boolean FALSE instead of a blank text placeholder.
#    B19 is a "slave" cell inside the merged range A19:G19, so a direct
#    Range.Value write on it is silently dropped by the merge logic.
#    Routing the write through Copy / PasteSpecial(values) lands the value
#    even on a merged slave cell, which matches how this file was produced.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Value = $false
$ws.Range("Z1").Copy()
$ws.Range("B19").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Fix the FLOOR() calls that were being called with an extra 3rd
#    argument (FLOOR only takes number + significance).
# ---------------------------------------------------------------------------
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
